$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (the "Förändrad" date column) for rows 2 through 12
# from 2023-10-08 (serial 45207) to 2023-10-09 (serial 45208)
for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    if ($current.Year -eq 2023 -and $current.Month -eq 10 -and $current.Day -eq 8) {
        $cell.Value = 45208
    }
}
